# Economic Dashboard update - 2025-11-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: GDP Nowcast date cell loses its "recently updated" yellow highlight ---
# Copy the (non-highlighted) date format from N7 (style index 47) onto C7
# (which currently uses the highlighted style index 48), preserving its value.
$ws.Range("N7").Copy()
$ws.Range("C7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 28 (Durable Orders, M/M % Delta) ---
$ws.Range("F28").Value = 0.0292251268148207
$ws.Range("G28").Value = -0.02746655290430811

# --- Row 29 (Durable Orders Y/Y % Delta | 5yr,5yr Fwd Inflation) ---
$ws.Range("G29").Value = 0.03397968857635882
$ws.Range("N29").Value = 45973
$ws.Range("Q29").Value = 2.17
$ws.Range("R29").Value = ""
$ws.Range("S29").Value = 2.2
$ws.Range("T29").Value = ""
$ws.Range("U29").Value = ""

# --- Row 30 (Dur Orders Non Def x Aircraft | 10yr TIPS) ---
$ws.Range("N30").Value = 45973
$ws.Range("Q30").Value = 2.27
$ws.Range("R30").Value = ""
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = ""
$ws.Range("U30").Value = ""

# --- Row 47 (FFR date only) ---
$ws.Range("N47").Value = 45972

# --- Row 48 (2y UST) ---
$ws.Range("N48").Value = 45971
$ws.Range("Q48").Value = 3.58
$ws.Range("R48").Value = ""
$ws.Range("S48").Value = ""
$ws.Range("T48").Value = 3.55
$ws.Range("U48").Value = 3.57

# --- Row 49 (5y UST) ---
$ws.Range("N49").Value = 45971
$ws.Range("Q49").Value = 3.72
$ws.Range("R49").Value = ""
$ws.Range("S49").Value = ""
$ws.Range("T49").Value = 3.67
$ws.Range("U49").Value = 3.69

# --- Row 50 (10y UST) ---
$ws.Range("N50").Value = 45971
$ws.Range("Q50").Value = 4.13
$ws.Range("R50").Value = ""
$ws.Range("S50").Value = ""
$ws.Range("T50").Value = 4.11
$ws.Range("U50").Value = 4.11

# --- Row 51: 30y Mortgage date cell loses its "recently updated" yellow highlight ---
$ws.Range("N7").Copy()
$ws.Range("N51").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 52 (BAA) ---
$ws.Range("N52").Value = 45971
$ws.Range("R52").Value = ""
$ws.Range("S52").Value = ""
$ws.Range("T52").Value = 5.86
$ws.Range("U52").Value = 5.83
